$ws = $excel.ActiveWorkbook.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-09T07:01:32.302554+00:00"
$ws.Range("K3").Value = "2025-11-09T07:01:32.302592+00:00"
$ws.Range("K4").Value = "2025-11-09T07:01:32.302614+00:00"
$ws.Range("K5").Value = "2025-11-09T07:01:36.291526+00:00"
$ws.Range("K6").Value = "2025-11-09T07:01:36.291555+00:00"
$ws.Range("K7").Value = "2025-11-09T07:01:36.291574+00:00"
$ws.Range("K8").Value = "2025-11-09T07:01:38.752747+00:00"
$ws.Range("K9").Value = "2025-11-09T07:01:41.480364+00:00"
$ws.Range("K10").Value = "2025-11-09T07:01:41.480395+00:00"
$ws.Range("K11").Value = "2025-11-09T07:01:41.480414+00:00"
$ws.Range("K12").Value = "2025-11-09T07:01:44.328725+00:00"
$ws.Range("K13").Value = "2025-11-09T07:01:44.328756+00:00"
$ws.Range("K14").Value = "2025-11-09T07:01:44.328775+00:00"
$ws.Range("K15").Value = "2025-11-09T07:01:44.328791+00:00"
$ws.Range("K16").Value = "2025-11-09T07:01:49.981378+00:00"
$ws.Range("K17").Value = "2025-11-09T07:01:52.290199+00:00"
$ws.Range("K18").Value = "2025-11-09T07:01:54.669928+00:00"
$ws.Range("K19").Value = "2025-11-09T07:01:54.669958+00:00"
$ws.Range("K20").Value = "2025-11-09T07:01:54.669976+00:00"
$ws.Range("K21").Value = "2025-11-09T07:01:57.414117+00:00"
$ws.Range("K22").Value = "2025-11-09T07:01:59.690018+00:00"
$ws.Range("K23").Value = "2025-11-09T07:01:59.690057+00:00"
$ws.Range("K24").Value = "2025-11-09T07:02:02.436000+00:00"
$ws.Range("K25").Value = "2025-11-09T07:02:02.436032+00:00"
$ws.Range("K26").Value = "2025-11-09T07:02:02.436051+00:00"
$ws.Range("K27").Value = "2025-11-09T07:02:04.699123+00:00"
$ws.Range("K28").Value = "2025-11-09T07:02:04.699151+00:00"
$ws.Range("K29").Value = "2025-11-09T07:02:04.699169+00:00"
$ws.Range("K30").Value = "2025-11-09T07:02:04.699195+00:00"
$ws.Range("K31").Value = "2025-11-09T07:02:04.699210+00:00"
$ws.Range("K32").Value = "2025-11-09T07:02:07.576189+00:00"
$ws.Range("K33").Value = "2025-11-09T07:02:07.576216+00:00"
$ws.Range("K34").Value = "2025-11-09T07:02:09.938292+00:00"
$ws.Range("K35").Value = "2025-11-09T07:02:09.938320+00:00"
$ws.Range("K36").Value = "2025-11-09T07:02:09.938337+00:00"
$ws.Range("K37").Value = "2025-11-09T07:02:12.662184+00:00"
$ws.Range("K38").Value = "2025-11-09T07:02:12.662221+00:00"
$ws.Range("K39").Value = "2025-11-09T07:02:12.662239+00:00"
$ws.Range("K40").Value = "2025-11-09T07:02:14.954024+00:00"
$ws.Range("K41").Value = "2025-11-09T07:02:14.954053+00:00"
$ws.Range("K42").Value = "2025-11-09T07:02:14.954071+00:00"
$ws.Range("K43").Value = "2025-11-09T07:02:14.954088+00:00"
$ws.Range("K44").Value = "2025-11-09T07:02:14.954103+00:00"
$ws.Range("K45").Value = "2025-11-09T07:02:14.954119+00:00"
$ws.Range("K46").Value = "2025-11-09T07:02:17.721687+00:00"
$ws.Range("K47").Value = "2025-11-09T07:02:17.721722+00:00"
$ws.Range("K48").Value = "2025-11-09T07:02:22.844268+00:00"
$ws.Range("K49").Value = "2025-11-09T07:02:22.844298+00:00"
$ws.Range("K50").Value = "2025-11-09T07:02:22.844316+00:00"
$ws.Range("K51").Value = "2025-11-09T07:02:25.636361+00:00"
$ws.Range("K52").Value = "2025-11-09T07:02:25.636392+00:00"
